$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.275.21"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.618.76"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.00"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.78"
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.24"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.618.93"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.286.81"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.30"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.57"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.33"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.90"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.95"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.16"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0525"
$ws.Range("E30").Value = "  +10.60%  "
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.42"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.179.08"
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.804"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.787"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("E43").Value = "  +4.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.756.68"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.94"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  +14.40%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.78"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -0.24%  "
